$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.001754667048134761
$ws.Range("C2").Value = 87981.0709163148
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 159515.9479795073
